$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells for the "mittlere" (averaged) columns ---
$ws.Range("M1").Value = "mtauf"
$ws.Range("N1").Value = "mtab"
$ws.Range("M1:N1").NumberFormat = "0.0000"

# Style the rest of the M/N columns (same numeric style as B/C) so every
# row in the used range gets a formatted, if empty, cell.
$ws.Range("M2:N31").NumberFormat = "0.00"

# --- Row 2: stand-alone averages (first block, not part of a fill-down) ---
$ws.Range("M2").Formula = "=AVERAGE(B2:B4)"
$ws.Range("N2").Formula = "=AVERAGE(C2:C4)"

# Rows 3 and 4 belong to the first triplet (already averaged in row 2)
# and stay blank, just carrying the number format.
$ws.Range("M3:M4").NumberFormat = "0.00"
$ws.Range("N3:N4").NumberFormat = "0.00"

# --- Rows 5:31: fill down the averages for every group of three
# measurements. This recreates the shared-formula block Excel builds
# when you drag-fill a formula down a column. ---
$ws.Range("M5:M31").Formula = "=AVERAGE(B5:B7)"
$ws.Range("N5:N31").Formula = "=AVERAGE(C5:C7)"

# Only every third row (the first row of each triplet) is meant to keep
# the formula - the two rows below each measurement group are blanked
# back out again, exactly like the existing E/F -> G/H/I pattern.
$ws.Range("M6:M7").ClearContents()
$ws.Range("N6:N7").ClearContents()

$ws.Range("M9:M10").ClearContents()
$ws.Range("N9:N10").ClearContents()

$ws.Range("M12:M13").ClearContents()
$ws.Range("N12:N13").ClearContents()

$ws.Range("M15:M16").ClearContents()
$ws.Range("N15:N16").ClearContents()

$ws.Range("M18:M19").ClearContents()
$ws.Range("N18:N19").ClearContents()

$ws.Range("M21:M22").ClearContents()
$ws.Range("N21:N22").ClearContents()

$ws.Range("M24:M25").ClearContents()
$ws.Range("N24:N25").ClearContents()

$ws.Range("M27:M28").ClearContents()
$ws.Range("N27:N28").ClearContents()

$ws.Range("M30:M31").ClearContents()
$ws.Range("N30:N31").ClearContents()

# Row 26 was re-typed by hand afterwards (small second-guess mentioned in
# the commit message), so it no longer shares the fill-down formula.
$ws.Range("M26").Formula = "=AVERAGE(B26:B28)"

# --- View state: scroll down a bit and leave the selection on N3, the
# last cell touched while filling in the new columns. ---
$ws.Range("N3").Select()
$excel.ActiveWindow.ScrollRow = 6
